# "Added last minute updates"
#
# The first paragraph of the document (the hidden **ID__...__ID** merge
# marker) is updated to:
#   1. Gain the same zero-width "padding" paragraph border (5-twip space
#      on all four sides, no visible rule) that every other body
#      paragraph in this template already carries.
#   2. Have its left indent widened from 120 to 225 twips (6pt -> 11.25pt)
#      to line up with the other bordered paragraphs.
#   3. Have its merge-field id text corrected from the stale
#      "AFFARS_pgi_5304_topic_15" id to "AFFARS_SMC_PGI_5304_606", and
#      lose the trailing run that held nothing but a single space
#      (collapsing the paragraph down to one run).

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# --- 1 & 2: paragraph formatting -----------------------------------
$b = $p.Borders
$b.DistanceFromTop = 5
$b.DistanceFromBottom = 5
$b.DistanceFromLeft = 5
$b.DistanceFromRight = 5

$p.LeftIndent = 11.25

# --- 3: text content --------------------------------------------------
# Grab the paragraph's range but exclude the trailing paragraph mark, so
# we replace exactly "**ID__AFFARS_pgi_5304_topic_15__ID** " (36 chars of
# id text + the lone space that used to live in its own run) with the
# new id text and nothing else. Word collapses the result back down to a
# single run.
$rng = $p.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "**ID__AFFARS_SMC_PGI_5304_606__ID**"
